$d = $word.ActiveDocument

# Locate the "Author" styled paragraph whose text is exactly "Edison Achalma"
# (the author byline right under the title heading).
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($p.Style.NameLocal -eq "Author" -and ($t -eq "Edison Achalma" -or $t -eq "Edison Achalma`r")) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Edison Achalma' Author paragraph"
}

# Word's InsertParagraphAfter (as exposed here) splits a range at its
# *start*, so to land the new, empty paragraph immediately after the
# target we instead split at the *start* of the paragraph that follows it.
$following = $target.Next()
$insertionPoint = $following.Range.Duplicate
$insertionPoint.Collapse(1)
$insertionPoint.InsertParagraphAfter()

# The paragraph that is now right after $target is the freshly created one;
# give it the "Author" style and the new affiliation text.
$newPara = $target.Next()
$newPara.Style = "Author"
$newPara.Range.Text = "Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga"
